## Common_MetaData.xlsx update
## - Adds a new "Links" sheet entry (row 28): BPPDIGITALRFILINK -> RFI link
##   used by the newly created RFI-forms test.
## - Keeps the existing selection on the "Links" sheet pointing at the
##   newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Links")

# --- Values for the new row -------------------------------------------------
# Shared strings must be introduced in this order so the resulting
# sharedStrings.xml lists the URL before the constant name, matching the
# order the strings were authored in.
$ws.Range("B28").Value = "https://web-stage-bppdigital.bppuniversity.com/contact-bpp/request-information"
$ws.Range("A28").Value = "BPPDIGITALRFILINK"
$ws.Range("C28").Value = "descr."

# --- Formatting matching the rest of the "Links" table ----------------------
$ws.Range("A28").HorizontalAlignment = -4108   # xlCenter, like A2:A27
$ws.Range("C28").HorizontalAlignment = -4108   # xlCenter, like C2:C27

# --- Hyperlink on the url cell ----------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B28"), "https://web-stage-bppdigital.bppuniversity.com/contact-bpp/request-information")

# the link column keeps plain (non hyperlink-blue) text, same as B22:B27
$ws.Range("B28").Font.Underline = -4142        # xlUnderlineStyleNone
$ws.Range("B28").Font.ThemeColor = 1           # xlThemeColorDark1 (automatic/black)
$ws.Range("B28").Font.TintAndShade = 0

# --- Selection / view state ---------------------------------------------------
$ws.Activate()
$ws.Range("A28").Select()

Write-Host "Added BPPDIGITALRFILINK row to Links sheet"
